$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Förändrad" (column C) date values for rows 2-6 from 45174 (2023-09-05)
# to 45175 (2023-09-06), matching the automatic update of the source data.
$ws.Range("C2:C6").Value = 45175
